# Rename metadata sheet 'General' to 'Table' and make it the active/selected
# sheet (it was the first sheet, now becomes the one shown on open).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("General")
$ws.Name = "Table"
$ws.Activate()
